$wb = $excel.ActiveWorkbook

# Rename the "Geno" worksheet to "Patient details".
# Renaming a sheet automatically updates references that point at it
# (e.g. the _xlnm._FilterDatabase defined name / AutoFilter range).
$ws = $wb.Worksheets.Item("Geno")
$ws.Name = "Patient details"

# Move the active selection from G10 to E9. Selecting this cell also
# resets the sheet's scrolled view, dropping the previous
# topLeftCell="B1" pin so column A is visible again.
$ws.Activate()
$ws.Range("E9").Select()
